# Rename the inline picture shapes' internal `name` metadata (the
# wp:docPr/@name and pic:cNvPr/@name attributes baked into each picture's
# drawing XML) for the Pearson logo (footers) and the BTEC logo (headers).
#
# The Word object model's InlineShape.Name setter only patches the
# wp:docPr/@name half of this pair, so we round-trip the shape's own
# Range.WordOpenXML (a standard, supported Word COM mechanism for reading/
# writing the raw WordprocessingML behind a range) and do a scoped literal
# replace of the old `name="..."` value with the new one. Because the
# snippet returned for a single picture's range only ever contains that
# one picture's name twice (docPr + cNvPr), a plain string replace is safe
# and touches nothing else (ids, relationships, extents are untouched).

function Rename-InlineShapeXmlName {
    param($Shape, $OldName, $NewName)

    $range = $Shape.Range
    $xml = $range.WordOpenXML
    $needle = 'name="' + $OldName + '"'
    $replacement = 'name="' + $NewName + '"'
    $updated = $xml.Replace($needle, $replacement)
    $range.WordOpenXML = $updated
}

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-InlineShapeXmlName $shp "image1.jpg" "image2.jpg"
                }
            }
        }

        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    Rename-InlineShapeXmlName $shp "image2.png" "image1.png"
                }
            }
        }
    }
}
